$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('COLEGIO E CURSO MATRIZ EDUCACAO BANGU', 1877, 1219, 1069, 379, 877),
    @('COLEGIO E CURSO MATRIZ EDUCACAO CAMPO GRANDE', 1889, 1564, 1462, 555, 845),
    @('COLÉGIO LEONARDO DA VINCI - ALFA', 1220, 717, 519, 371, 834),
    @('APOGEU BENFICA', 1142, 873, 800, 295, 790),
    @('APOGEU SANTO ANTÔNIO II', 1608, 973, 858, 377, 751),
    @('COLÉGIO LEONARDO DA VINCI - BETA', 715, 438, 327, 238, 749),
    @('APOGEU GLOBAL SCHOOL FERREIRA GUIMARÃES', 1010, 903, 514, 374, 733),
    @('APOGEU SANTO ANTÔNIO I', 1943, 1585, 1481, 518, 707),
    @('COLEGIO E CURSO MATRIZ EDUCACAO ROCHA MIRANDA', 765, 580, 526, 175, 700),
    @('APOGEU GLOBAL SCHOOL CIDADE ALTA', 388, 189, 121, 69, 673),
    @('COLEGIO E CURSO MATRIZ EDUCACAO TAQUARA', 1196, 883, 785, 284, 618),
    @('ESCOLA SA PEREIRA (R. DA MATRIZ)', 1090, 453, 268, 163, 593),
    @('ESCOLA SAP', 1838, 375, 200, 132, 566),
    @('COLEGIO QI METROPOLITANO', 1079, 499, 232, 203, 566),
    @('COLÉGIO QI RIO 2', 1171, 415, 323, 279, 558),
    @('COLEGIO SARAH DAWSEY - JUIZ DE FORA', 1152, 447, 284, 203, 557),
    @('COLÉGIO E CURSO CUBO BARRA GOLF', 2289, 959, 675, 464, 536),
    @('COLEGIO E CURSO MATRIZ EDUCACAO RETIRO DOS ARTISTAS', 960, 716, 646, 354, 524),
    @('COLÉGIO QI TIJUCA', 1042, 224, 169, 142, 503),
    @('COLÉGIO AMERICANO', 481, 290, 224, 164, 481),
    @('COLEGIO E CURSO MATRIZ EDUCACAO MADUREIRA', 979, 739, 676, 217, 464),
    @('UNIFICADO ZONA SUL', 1254, 577, 393, 220, 447),
    @('COLEGIO E CURSO MATRIZ EDUCACAO SÃO JOÃO DE MERITI', 977, 617, 550, 152, 446),
    @('COLÉGIOS INTEGRADOS LEONARDO DA VINCI - GAMA', 815, 369, 265, 191, 443),
    @('COLÉGIO E CURSO AO CUBO BOTAFOGO', 1212, 386, 232, 164, 366),
    @('COLEGIO E CURSO MATRIZ EDUCACAO DUQUE DE CAXIAS', 918, 674, 618, 287, 347),
    @('COLEGIO QI RECREIO', 866, 298, 231, 213, 334),
    @('COLEGIO E CURSO MATRIZ EDUCACAO NOVA IGUACU', 1069, 811, 737, 179, 321),
    @('COLÉGIO QI FREGUESIA', 800, 216, 126, 105, 321),
    @('COLÉGIO QI VALQUEIRE', 2693, 765, 610, 563, 292),
    @('CRECHE ESCOLA GLOBAL TREE - RIO 2', 959, 361, 270, 189, 248),
    @('COLEGIO E CURSO UNIAO LTDA', 242, 122, 95, 75, 204),
    @('COLÉGIO E CURSO MATRIZ EDUCACAO TIJUCA', 1208, 625, 544, 120, 188),
    @('CRECHE ESCOLA GLOBAL TREE - PENINSULA', 347, 135, 82, 68, 182),
    @('CRECHE ESCOLA GLOBAL TREE - ABM', 598, 290, 229, 189, 159),
    @('CRECHE ESCOLA GLOBAL TREE - RECREIO', 1090, 397, 237, 187, 136),
    @('APOGEU - PARÁ DE MINAS', 1, 0, 0, 0, 0),
    @('COLÉGIO UNIFICADO INDEPENDÊNCIA', 1, 0, 0, 0, 0),
    @('Leads Raiz Sem Unidades', 4020, 913, 320, 143, 0)
)

# Clear the rows that will no longer be used (old data went to row 44)
$ws.Range("A2:F44").ClearContents()

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value2 = $row[$j]
    }
}

Write-Host "Rows written: $($data.Count)"
Write-Host $ws.UsedRange.Address()
